$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:87 down to 44:88.
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the new data record.
$ws.Range("A43").Value2 = 1
$ws.Range("B43").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value2 = "Arica y Parinacota"
$ws.Range("D43").Value2 = 45118
$ws.Range("E43").Value2 = 15
$ws.Range("F43").Value2 = 100112031
$ws.Range("G43").Value2 = "Poroto verde"
$ws.Range("H43").Value2 = "Magnum"
$ws.Range("I43").Value2 = "Primera"
$ws.Range("J43").Value2 = 400
$ws.Range("K43").Value2 = 17000
$ws.Range("L43").Value2 = 18000
$ws.Range("M43").Value2 = 17500
$ws.Range("N43").Value2 = "`$/malla 25 kilos"
$ws.Range("O43").Value2 = "Perú"
$ws.Range("P43").Value2 = 700
$ws.Range("Q43").Value2 = 25
$ws.Range("R43").Value2 = "Hortaliza"
